$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.548.72'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').Value = '1.884.15'
$ws.Range('E3').Value = '  +1.50%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''248.25'
$ws.Range('E5').Value = '  +6.72%  '
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = '''0.4768'
$ws.Range('E7').Value = '  +1.58%  '
$ws.Range('D8').Value = '''0.2923'
$ws.Range('E8').Value = '  +3.52%  '
$ws.Range('D9').Value = '''0.06536'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('D10').Value = '''22.10'
$ws.Range('E10').Value = '  +6.95%  '
$ws.Range('B11').Value = 'Litecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D11').Value = '''97.97'
$ws.Range('E11').Value = '  +5.31%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '''0.07726'
$ws.Range('E12').Value = '  +0.85%  '
$ws.Range('D13').Value = '''0.7416'
$ws.Range('E13').Value = '  +9.88%  '
$ws.Range('D14').Value = '1.885.99'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').Value = '''5.156'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').Value = '''274.59'
$ws.Range('E16').Value = '  +3.32%  '
$ws.Range('D17').Value = '30.554.69'
$ws.Range('D18').Value = '''13.51'
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('D19').Value = '''0.000007575'
$ws.Range('E19').Value = '  +1.48%  '
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').Value = '2.132.93'
$ws.Range('E21').Value = '  +1.94%  '
$ws.Range('D22').Value = '''0.9999'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '''5.264'
$ws.Range('E23').Value = '  +2.83%  '
$ws.Range('D24').Value = '''6.208'
$ws.Range('E24').Value = '  +2.49%  '
$ws.Range('D25').Value = '''9.320'
$ws.Range('E25').Value = '  +1.01%  '
$ws.Range('D26').Value = '''163.45'
$ws.Range('E26').Value = '  -1.03%  '
$ws.Range('D27').Value = '''18.95'
$ws.Range('E27').Value = '  +2.92%  '
$ws.Range('D28').Value = '''1.948'
$ws.Range('E28').Value = '  +4.23%  '
$ws.Range('D29').Value = '''0.1007'
$ws.Range('E29').Value = '  +3.05%  '
$ws.Range('D30').Value = '''1.367'
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = '''1.518'
$ws.Range('E31').Value = '  +4.94%  '
$ws.Range('D32').Value = '''4.336'
$ws.Range('E32').Value = '  +3.94%  '
$ws.Range('D33').Value = '''4.120'
$ws.Range('E33').Value = '  +4.08%  '
$ws.Range('D34').Value = '''0.04827'
$ws.Range('E34').Value = '  +4.55%  '
$ws.Range('E35').Value = '  +2.47%  '
$ws.Range('D36').Value = '''0.7032'
$ws.Range('E36').Value = '  +3.34%  '
$ws.Range('E38').Value = '  +3.72%  '
$ws.Range('D39').Value = '''2.754'
$ws.Range('E39').Value = '  +1.58%  '
$ws.Range('D40').Value = '''6.333'
$ws.Range('E40').Value = '  +1.32%  '
$ws.Range('D41').Value = '''1.999'
$ws.Range('E41').Value = '  +7.31%  '
$ws.Range('D42').Value = '''71.80'
$ws.Range('E42').Value = '  +2.78%  '
$ws.Range('D43').Value = '''0.4229'
$ws.Range('E43').Value = '  +5.20%  '
$ws.Range('D44').Value = '''0.8435'
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').Value = '''102.95'
$ws.Range('E46').Value = '  +0.73%  '
$ws.Range('D47').Value = '''9.338'
$ws.Range('E47').Value = '  +2.78%  '
$ws.Range('D48').Value = '''7.110'
$ws.Range('E48').Value = '  +3.74%  '
$ws.Range('E49').Value = '  +4.82%  '
$ws.Range('D50').Value = '''920.05'
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('E51').Value = '  +4.84%  '
